$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ReqPow_AC)
$ws.Range("B2").Value = 98.04961276102425
$ws.Range("C2").Value = 196.0992255220485
$ws.Range("D2").Value = 1258.85843638851
$ws.Range("E2").Value = 1960.992255220485
$ws.Range("F2").Value = 1960.992255220485
$ws.Range("G2").Value = 1598.469839612929
$ws.Range("H2").Value = 1598.469839612929
$ws.Range("I2").Value = 1494.6688982533
$ws.Range("J2").Value = 1494.6688982533
$ws.Range("K2").Value = 196.0992255220485
$ws.Range("L2").Value = 1960.992255220485
$ws.Range("M2").Value = 1960.992255220485
$ws.Range("N2").Value = 1494.6688982533
$ws.Range("O2").Value = 1494.6688982533
$ws.Range("P2").Value = 1494.6688982533
$ws.Range("Q2").Value = 1494.6688982533
$ws.Range("R2").Value = 196.0992255220485
$ws.Range("S2").Value = 196.0992255220485
$ws.Range("T2").Value = 196.0992255220485
$ws.Range("U2").Value = 98.04961276102425

# Row 3 (ReqPow_FC)
$ws.Range("B3").Value = 98.04961276102425
$ws.Range("C3").Value = 196.0992255220485
$ws.Range("D3").Value = 1005.334193409084
$ws.Range("E3").Value = 1580.705890751345
$ws.Range("F3").Value = 1580.705890751345
$ws.Range("G3").Value = 1655.735457629397
$ws.Range("H3").Value = 1655.735457629397
$ws.Range("I3").Value = 1494.6688982533
$ws.Range("J3").Value = 1494.6688982533
$ws.Range("K3").Value = 196.0992255220485
$ws.Range("L3").Value = 1580.705890751345
$ws.Range("M3").Value = 1580.705890751345
$ws.Range("N3").Value = 1494.6688982533
$ws.Range("O3").Value = 1494.6688982533
$ws.Range("P3").Value = 1494.6688982533
$ws.Range("Q3").Value = 1494.6688982533
$ws.Range("R3").Value = 196.0992255220485
$ws.Range("S3").Value = 196.0992255220485
$ws.Range("T3").Value = 196.0992255220485
$ws.Range("U3").Value = 98.04961276102425

# Row 4 (ReqPow_Batt)
$ws.Range("D4").Value = 253.5242429794263
$ws.Range("E4").Value = 380.2863644691395
$ws.Range("F4").Value = 380.2863644691395
$ws.Range("G4").Value = -57.26561801646823
$ws.Range("H4").Value = -57.26561801646823
$ws.Range("L4").Value = 380.2863644691395
$ws.Range("M4").Value = 380.2863644691395
